$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 718 entirely ("「お金で買えぬ10のもの」" post) — everything below shifts up by one.
$ws.Rows.Item(718).Delete()
